$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price/volume cells to match the latest scrape.
# Leading "'" forces text so numeric-looking strings (prices, percents)
# keep exact formatting (trailing zeros, thousands separators); Style is
# reset to Normal afterwards so no stray number-format style is left behind.
$ws.Range("D2").Value = "'63.304.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.09%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.617.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.28%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.43%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'605.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.86%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'146.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.71%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.38%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.109"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.66%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'5.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.54%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.372"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.66%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.35%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'27.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.35%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.101.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.66%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'63.196.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.02%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000147"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.34%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.673.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.78%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'11.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.90%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.72%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'342.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.70%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.00%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.92%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'66.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.54%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.43%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +7.11%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'558.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.162"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.61%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Binance-PegBSC-USD"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.22%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Aptos"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'7.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.77%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0₃0847"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.28%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.84%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.70%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'167.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.25%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'FirstDigitalUSD"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.10%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'PolygonEcosystemToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.403"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +6.51%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'19.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.40%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.11%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'165.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.33%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.33%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'21.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.85%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0567"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.625"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.93%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0246"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.57%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0958"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +11.95%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'18.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.19%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.181"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.02%  "
$ws.Range("E51").Style = "Normal"
